# Updated capital structure database
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.11445
$ws.Range("E2").Value = 0.07730000000000001
$ws.Range("F2").Value = 0.025
$ws.Range("G2").Value = 0.1469165805392611
$ws.Range("H2").Value = 0.1469165805392611
$ws.Range("I2").Value = 0.1378378659984185
$ws.Range("J2").Value = 0.1289460748298135
$ws.Range("K2").Value = 27175.1
$ws.Range("L2").Value = 0.0792725100282638
$ws.Range("M2").Value = 16030.4
$ws.Range("N2").Value = 0.04010196602299989
$ws.Range("O2").Value = 0.5898929534757921
$ws.Range("P2").Value = 15296.6
$ws.Range("Q2").Value = 0.03826627741462597
$ws.Range("R2").Value = 0.5628902929520039
$ws.Range("S2").Value = 733.7999999999993
$ws.Range("T2").Value = 0.04577552649965062
$ws.Range("U2").Value = 77360.3
$ws.Range("V2").Value = 0.1935260581226344
$ws.Range("W2").Value = 0.1480299425380022
$ws.Range("X2").Value = 0.06862902734486423
$ws.Range("Y2").Value = 0.07940091519313799
$ws.Range("Z2").Value = 0.9700369365408266
$ws.Range("AA2").Value = 0.148660898225385
$ws.Range("AB2").Value = 0.06093873474713649
$ws.Range("AC2").Value = 0.08772216347824849
$ws.Range("AD2").Value = 329922.7
$ws.Range("AE2").Value = 14.29362379766415
$ws.Range("AF2").Value = 329936.9936237977
$ws.Range("AG2").Value = 252576.6936237977
$ws.Range("AH2").Value = 0.4521679377847652
$ws.Range("AI2").Value = 0.6009916022815114
$ws.Range("AJ2").Value = 0.3871989003098585
$ws.Range("AK2").Value = 0.5355428820262085
$ws.Range("AL2").Value = 12097.1
$ws.Range("AM2").Value = 12097.1
$ws.Range("AN2").Value = 6.736550592526557
$ws.Range("AO2").Value = 3.90564680791264
$ws.Range("AP2").Value = 5.157255548314176
$ws.Range("AQ2").Value = 3.90564680791264

# Row 3 (E3 cell is cleared entirely, D3 updated, no F3)
$ws.Range("D3").Value = 1.147
$ws.Range("E3").ClearContents()
$ws.Range("G3").Value = 0.08215619418736808
$ws.Range("H3").Value = 0.08215619418736808
$ws.Range("I3").Value = 0.05527866134769722
$ws.Range("J3").Value = 0.05342762550348209
$ws.Range("K3").Value = 196.7
$ws.Range("L3").Value = 0.03193700276018834
$ws.Range("M3").Value = 97.2
$ws.Range("N3").Value = 0.02637362637362637
$ws.Range("O3").Value = 0.4941535332994408
$ws.Range("P3").Value = 97.2
$ws.Range("Q3").Value = 0.02637362637362637
$ws.Range("R3").Value = 0.4941535332994408
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0
$ws.Range("U3").Value = 3381.4
$ws.Range("V3").Value = 0.9174874508207842
$ws.Range("W3").Value = 0.06916315049226442
$ws.Range("X3").Value = 0.06463553811180858
$ws.Range("Y3").Value = 0.004527612380455834
$ws.Range("Z3").Value = 8.290554397970546
$ws.Range("AA3").Value = 0.4429446355910167
$ws.Range("AB3").Value = 0.06115882700793634
$ws.Range("AC3").Value = 0.3817858085830804
$ws.Range("AD3").Value = 290.1
$ws.Range("AE3").Value = 14.29362379766415
$ws.Range("AF3").Value = 304.3936237976642
$ws.Range("AG3").Value = -3077.006376202336
$ws.Range("AH3").Value = 0.07629116274732557
$ws.Range("AI3").Value = 0.05734950910370558
$ws.Range("AJ3").Value = -5.056760261510153
$ws.Range("AK3").Value = -1.597371417414577
$ws.Range("AN3").Value = 0.8206042090970809
$ws.Range("AP3").Value = -8.703910319649061

# Row 4
$ws.Range("D4").Value = 0.0475
$ws.Range("E4").Value = 0.07730000000000001
$ws.Range("F4").Value = 0.0194
$ws.Range("G4").Value = 0.06919331729865348
$ws.Range("H4").Value = 0.06919331729865348
$ws.Range("I4").Value = 0.07238944822981024
$ws.Range("J4").Value = 0.0709549012374512
$ws.Range("K4").Value = 1864.4
$ws.Range("L4").Value = 0.06043475160292774
$ws.Range("M4").Value = 683
$ws.Range("N4").Value = 0.03028556225611919
$ws.Range("O4").Value = 0.3663376957734392
$ws.Range("P4").Value = 683
$ws.Range("Q4").Value = 0.03028556225611919
$ws.Range("R4").Value = 0.3663376957734392
$ws.Range("U4").Value = 1467.1
$ws.Range("V4").Value = 0.0650540971975878
$ws.Range("W4").Value = 0.167005562671874
$ws.Range("X4").Value = 0.07053245931422591
$ws.Range("Y4").Value = 0.09647310335764811
$ws.Range("Z4").Value = 2.121500532957398
$ws.Range("AA4").Value = 0.1505308607911922
$ws.Range("AB4").Value = 0.06079011905119952
$ws.Range("AC4").Value = 0.08974074173999272
$ws.Range("AD4").Value = 5265.5
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 5265.5
$ws.Range("AG4").Value = 3798.4
$ws.Range("AH4").Value = 0.1892873191336389
$ws.Range("AI4").Value = 0.27590452985407
$ws.Range("AJ4").Value = 0.1441496144271055
$ws.Range("AK4").Value = 0.2156050268484566
$ws.Range("AL4").Value = 148.4
$ws.Range("AM4").Value = 148.4
$ws.Range("AN4").Value = 2.250213675213675
$ws.Range("AO4").Value = 15.04851752021563
$ws.Range("AP4").Value = 1.623247863247863
$ws.Range("AQ4").Value = 15.04851752021563

# Row 5
$ws.Range("D5").Value = 0.0919
$ws.Range("E5").Value = 0.0437
$ws.Range("F5").Value = 0.623
$ws.Range("G5").Value = 0.09504705337145039
$ws.Range("H5").Value = 0.09504705337145039
$ws.Range("I5").Value = 0.0935306134979129
$ws.Range("J5").Value = 0.08813303237953879
$ws.Range("K5").Value = 7018.8
$ws.Range("L5").Value = 0.05903265770596206
$ws.Range("M5").Value = 3787
$ws.Range("N5").Value = 0.02727213151115838
$ws.Range("O5").Value = 0.539550920385251
$ws.Range("P5").Value = 3787
$ws.Range("Q5").Value = 0.02727213151115838
$ws.Range("R5").Value = 0.539550920385251
$ws.Range("U5").Value = 8020.1
$ws.Range("V5").Value = 0.05775685818131538
$ws.Range("W5").Value = 0.1290543224041305
$ws.Range("X5").Value = 0.06672559537550253
$ws.Range("Y5").Value = 0.06232872702862792
$ws.Range("Z5").Value = 1.665560933242745
$ws.Range("AA5").Value = 0.1467909356595777
$ws.Range("AB5").Value = 0.06108735044307346
$ws.Range("AC5").Value = 0.08570358521650426
$ws.Range("AD5").Value = 18895
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 18895
$ws.Range("AG5").Value = 10874.9
$ws.Range("AH5").Value = 0.1197745613918317
$ws.Range("AI5").Value = 0.2322845390328997
$ws.Range("AJ5").Value = 0.07262783618482301
$ws.Range("AK5").Value = 0.1483127648344814
$ws.Range("AL5").Value = 2305.2
$ws.Range("AM5").Value = 2305.2
$ws.Range("AN5").Value = 1.630144077301354
$ws.Range("AO5").Value = 4.824093354155822
$ws.Range("AP5").Value = 0.9382193080838581
$ws.Range("AQ5").Value = 4.824093354155822

# Row 6
$ws.Range("D6").Value = 0.137
$ws.Range("E6").Value = 0.171
$ws.Range("F6").Value = 0.025
$ws.Range("G6").Value = 0.1948765224686518
$ws.Range("H6").Value = 0.1948765224686518
$ws.Range("I6").Value = 0.1795475023060411
$ws.Range("J6").Value = 0.1531496418830021
$ws.Range("K6").Value = 18095.2
$ws.Range("L6").Value = 0.09681734228498175
$ws.Range("M6").Value = 11463.2
$ws.Range("N6").Value = 0.04885362408893821
$ws.Range("O6").Value = 0.6334939652504531
$ws.Range("P6").Value = 10729.4
$ws.Range("Q6").Value = 0.04572633071915815
$ws.Range("R6").Value = 0.5929417746142623
$ws.Range("S6").Value = 733.7999999999993
$ws.Range("T6").Value = 0.06401353897689993
$ws.Range("U6").Value = 64491.7
$ws.Range("V6").Value = 0.2748493674241553
$ws.Range("W6").Value = 0.2017234652129806
$ws.Range("X6").Value = 0.1122853155817656
$ws.Range("Y6").Value = 0.08943814963121502
$ws.Range("Z6").Value = 0.7007232167963258
$ws.Range("AA6").Value = 0.1073155097114626
$ws.Range("AB6").Value = 0.05956231122263762
$ws.Range("AC6").Value = 0.04775319848882495
$ws.Range("AD6").Value = 305472.1
$ws.Range("AE6").Value = 0
$ws.Range("AF6").Value = 305472.1
$ws.Range("AG6").Value = 240980.4
$ws.Range("AH6").Value = 0.5655676864909921
$ws.Range("AI6").Value = 0.689162332970033
$ws.Range("AJ6").Value = 0.5066613515460315
$ws.Range("AK6").Value = 0.6362357548165115
$ws.Range("AL6").Value = 9643.5
$ws.Range("AM6").Value = 9643.5
$ws.Range("AN6").Value = 8.805641313904383
$ws.Range("AO6").Value = 3.479805050033701
$ws.Range("AP6").Value = 6.946581917239589
$ws.Range("AQ6").Value = 3.479805050033701
